$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.325.58"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.01%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.559.62"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.12%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9993"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.30%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.0000"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "287.96"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.25%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3779"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.45%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3274"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.34"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -8.79%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.141"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.36%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07385"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.17%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9992"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.31%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.38"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.96%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.858"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.40%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.763"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.89%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.552.55"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.50%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001079"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.45%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06648"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.71%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "86.23"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.86%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.421"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.17%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.13%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.16"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.50%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.69"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.05%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.306.08"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.09%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.297"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -4.35%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.579"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "151.51"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.48%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.38"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.55%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.941"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.45%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "122.78"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.33%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.723.55"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.83%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.087"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.18%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.936"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.09%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.895"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -5.57%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.439"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.62%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08215"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.27%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02371"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.47%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06317"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.06%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.331"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.98%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2156"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.99%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.250"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.14%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.09"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.74%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6087"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.14%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.16%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.76"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.91%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5939"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.00%  "

$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.744"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.74%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.988"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.53%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "122.92"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.70%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.179"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.35%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07087"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.60%  "
